# Fixed a bug in WinAreaData
# Restores the correct row order/values for the reel-weighting table on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(901, 16, 15, 45, 60, 60)
    3  = @(401, 9, 48, 67, 75, 45)
    4  = @(701, 3, 90, 45, 97, 15)
    5  = @(801, 3, 67, 65, 52, 45)
    7  = @(1203, 3, 15, 15, 15, 15)
    8  = @(101, 9, 30, 15, 60, 15)
    9  = @(902, 1, 0, 0, 0, 0)
    10 = @(301, 6, 45, 30, 60, 45)
    11 = @(501, 9, 52, 30, 75, 45)
    12 = @(201, 9, 30, 15, 45, 30)
    13 = @(1001, 18, 30, 75, 60, 72)
    16 = @(3, 0, 3, 3, 3, 3)
    19 = @(2, 0, 2, 2, 2, 2)
    20 = @(1, 0, 2, 2, 2, 2)
    21 = @(502, 0, 4, 0, 0, 0)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($col = 1; $col -le 6; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 1]
    }
}
